$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.016.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.26%  "

$ws.Range("E9").Value = "  -2.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.94%  "

$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.861.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.635.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.528"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0745"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.013.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.81%  "

$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "190.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.132"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  -2.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.12%  "

$ws.Range("E30").Value = "  -1.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0483"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.87%  "

$ws.Range("E33").Value = "  -3.99%  "

$ws.Range("E34").Value = "  -1.50%  "

$ws.Range("E35").Value = "  -2.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.873"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.132.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("E38").Value = "  -0.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.525"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.26%  "

$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.89%  "

$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("E43").Value = "  -3.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.771.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.28%  "

$ws.Range("E45").Value = "  -0.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0528"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.61%  "

$ws.Range("E48").Value = "  +1.70%  "

$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.08%  "
